# Actualizacion "Nuevo formato 15 jun 2021": recalculo de estadisticas por docente
# en las tres hojas (1er Parcial, 2o Parcial, Final) del libro "Totales Docentes".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
# Fila 2: C2, D2, E2, F2, I2
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 7.84
$ws.Range("E2").Value = 94
$ws.Range("F2").Value = 92.16
$ws.Range("I2").Value = 8.300000000000001
# Fila 22: C22, D22, E22, F22
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 4.9
$ws.Range("E22").Value = 145
$ws.Range("F22").Value = 71.08
# Fila 31: C31, D31, E31, F31, I31
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 143
$ws.Range("F31").Value = 99.31
$ws.Range("I31").Value = 7.7

$ws = $wb.Worksheets.Item("2o Parcial")
# Fila 2: C2, D2, E2, F2, I2
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 7.84
$ws.Range("E2").Value = 94
$ws.Range("F2").Value = 92.16
$ws.Range("I2").Value = 7.8
# Fila 6: C6, D6, E6, F6, G6, H6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 93
$ws.Range("F6").Value = 78.81
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 21.19
# Fila 19: C19, D19, I19
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("I19").Value = 7.1
# Fila 22: C22, D22, E22, F22, G22, H22, I22
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = 7.35
$ws.Range("E22").Value = 173
$ws.Range("F22").Value = 84.8
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 10.29
$ws.Range("I22").Value = 7.9
# Fila 31: C31, D31, E31, F31, I31
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 143
$ws.Range("F31").Value = 99.31
$ws.Range("I31").Value = 7.4

$ws = $wb.Worksheets.Item("Final")
# Fila 2: C2, D2, E2, F2, I2
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 7.84
$ws.Range("E2").Value = 94
$ws.Range("F2").Value = 92.16
$ws.Range("I2").Value = 8
# Fila 3: I3
$ws.Range("I3").Value = 9.199999999999999
# Fila 4: I4
$ws.Range("I4").Value = 8.1
# Fila 6: I6
$ws.Range("I6").Value = 7.3
# Fila 7: I7
$ws.Range("I7").Value = 7.7
# Fila 8: I8
$ws.Range("I8").Value = 7.3
# Fila 9: I9
$ws.Range("I9").Value = 7.8
# Fila 10: I10
$ws.Range("I10").Value = 7
# Fila 11: I11
$ws.Range("I11").Value = 6.8
# Fila 12: I12
$ws.Range("I12").Value = 7
# Fila 13: I13
$ws.Range("I13").Value = 6.6
# Fila 14: I14
$ws.Range("I14").Value = 7.1
# Fila 15: I15
$ws.Range("I15").Value = 7.3
# Fila 16: I16
$ws.Range("I16").Value = 8.199999999999999
# Fila 17: I17
$ws.Range("I17").Value = 7.7
# Fila 18: I18
$ws.Range("I18").Value = 7.8
# Fila 19: E19, F19, G19, H19, I19
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 94.59
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 5.41
$ws.Range("I19").Value = 7.4
# Fila 20: E20, F20, G20, H20, I20
$ws.Range("E20").Value = 211
$ws.Range("F20").Value = 84.06
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 11.55
$ws.Range("I20").Value = 6.7
# Fila 21: I21
$ws.Range("I21").Value = 7.6
# Fila 22: C22, D22, E22, F22, G22, H22, I22
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 4.9
$ws.Range("E22").Value = 187
$ws.Range("F22").Value = 91.67
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 3.43
$ws.Range("I22").Value = 7.7
# Fila 24: I24
$ws.Range("I24").Value = 6.8
# Fila 25: I25
$ws.Range("I25").Value = 8.199999999999999
# Fila 26: I26
$ws.Range("I26").Value = 8.1
# Fila 27: I27
$ws.Range("I27").Value = 7.4
# Fila 28: I28
$ws.Range("I28").Value = 7.1
# Fila 29: I29
$ws.Range("I29").Value = 7.7
# Fila 30: I30
$ws.Range("I30").Value = 7
# Fila 31: C31, D31, E31, F31, I31
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 143
$ws.Range("F31").Value = 99.31
$ws.Range("I31").Value = 7.4
# Fila 32: I32
$ws.Range("I32").Value = 7.4
# Fila 33: I33
$ws.Range("I33").Value = 7
# Fila 34: I34
$ws.Range("I34").Value = 8.199999999999999
# Fila 36: I36
$ws.Range("I36").Value = 7.7
# Fila 37: I37
$ws.Range("I37").Value = 8.5
# Fila 38: I38
$ws.Range("I38").Value = 7
# Fila 39: I39
$ws.Range("I39").Value = 6.2
# Fila 41: I41
$ws.Range("I41").Value = 7.1
# Fila 42: I42
$ws.Range("I42").Value = 7
# Fila 44: I44
$ws.Range("I44").Value = 7.9
# Fila 45: I45
$ws.Range("I45").Value = 8.699999999999999
# Fila 46: I46
$ws.Range("I46").Value = 8.4
# Fila 47: I47
$ws.Range("I47").Value = 7.4
# Fila 48: I48
$ws.Range("I48").Value = 7.5
# Fila 49: I49
$ws.Range("I49").Value = 7.2
# Fila 50: I50
$ws.Range("I50").Value = 8.4
# Fila 51: I51
$ws.Range("I51").Value = 8.1
